$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 31
$ws1.Range("F5").Value = 5221
$ws1.Range("F6").Value = 5221
$ws1.Range("F7").Value = 149
$ws1.Range("F8").Value = 135
$ws1.Range("F9").Value = 527
$ws1.Range("F10").Value = 5
$ws1.Range("F11").Value = 1176
$ws1.Range("F12").Value = 739
$ws1.Range("F13").Value = 5143
$ws1.Range("F14").Value = 28
$ws1.Range("F15").Value = 71
$ws1.Range("F16").Value = 89
$ws1.Range("F17").Value = 279
$ws1.Range("F18").Value = 279
$ws1.Range("F19").Value = 249
$ws1.Range("F20").Value = 102
$ws1.Range("F21").Value = 250
$ws1.Range("F22").Value = 3882
$ws1.Range("F24").Value = 3799
$ws1.Range("F25").Value = 183
$ws1.Range("F26").Value = 178
$ws1.Range("F28").Value = 231
$ws1.Range("F29").Value = 244
$ws1.Range("F31").Value = 112
$ws1.Range("F32").Value = 111
$ws1.Range("F36").Value = 17
$ws1.Range("F37").Value = 6746
$ws1.Range("F38").Value = 1095
$ws1.Range("F39").Value = 516
$ws1.Range("F40").Value = 101
$ws1.Range("F43").Value = 1379
$ws1.Range("F44").Value = 172
$ws1.Range("F45").Value = 686
$ws1.Range("F47").Value = 2304
$ws1.Range("F50").Value = 780
$ws1.Range("F51").Value = 926

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 1
$ws2.Range("F25").Value = 812

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 216

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 216
$ws4.Range("F7").Value = 5222
$ws4.Range("F8").Value = 5222
$ws4.Range("F9").Value = 149
$ws4.Range("F10").Value = 135
$ws4.Range("F12").Value = 527
$ws4.Range("F13").Value = 1176
$ws4.Range("F14").Value = 739
$ws4.Range("F15").Value = 5144
$ws4.Range("F16").Value = 28
$ws4.Range("F17").Value = 71
$ws4.Range("F18").Value = 89
$ws4.Range("F19").Value = 279
$ws4.Range("F20").Value = 279
$ws4.Range("F21").Value = 249
$ws4.Range("F22").Value = 102
$ws4.Range("F23").Value = 250
$ws4.Range("F24").Value = 3882
$ws4.Range("F25").Value = 3799
$ws4.Range("F26").Value = 183
$ws4.Range("F27").Value = 178
$ws4.Range("F28").Value = 231
$ws4.Range("F29").Value = 244
$ws4.Range("F31").Value = 112
$ws4.Range("F32").Value = 111
$ws4.Range("F35").Value = 17
$ws4.Range("F37").Value = 6746
$ws4.Range("F38").Value = 1095
$ws4.Range("F39").Value = 516
$ws4.Range("F41").Value = 101
$ws4.Range("F44").Value = 1379
$ws4.Range("F45").Value = 172
$ws4.Range("F46").Value = 686
$ws4.Range("F47").Value = 2304
$ws4.Range("F49").Value = 780
$ws4.Range("F50").Value = 926
